# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# Rows 16-18 stay the same (CC / 1049937241 / NOREIDIS RODRIGUEZ MANCILLA / 1607-1609 / 24640 / 925000).
# Rows 19-25 become NOREIDIS RODRIGUEZ MANCILLA rows with periods 1701,1612,1611,1610,1609,1608,1607.
# Row 26 becomes the LEONARDO FABIO RODRIGUEZ ROYERO row with period 1609 / 28640 / 2892800.

$periodos = @("1701", "1612", "1611", "1610", "1609", "1608", "1607")

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 19 + $i
    $ws.Range("B$row").Value = "CC"
    $ws.Range("C$row").Value = "1049937241"
    $ws.Range("D$row").Value = "NOREIDIS RODRIGUEZ MANCILLA"
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = 24640
    $ws.Range("G$row").Value = 925000
}

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "77165479"
$ws.Range("D26").Value = "LEONARDO FABIO RODRIGUEZ ROYERO"
$ws.Range("E26").Value = "1609"
$ws.Range("F26").Value = 28640
$ws.Range("G26").Value = 2892800
